# Muuta ELY_Turvallisuusvastaavan oikeudet urakkakohtaisiksi.
# Column P on the "Oikeudet" sheet holds the rights of the
# "ELY turvallisuusvastaava" role. Remove the "*" (and "+") qualifiers
# from the rights so the role's access becomes row/project specific
# rather than blanket "R*"/"R*,W+".

$wb = $excel.ActiveWorkbook
$wsOikeudet = $wb.Worksheets.Item("Oikeudet")
$wsRoolit = $wb.Worksheets.Item("Roolit")

# Rows whose P-column value is exactly "R*" -> becomes "R"
$rRowsOnly = @(7,23,24,25,26,45,63,64,65,66,67,68,69,70,71,72,77,78,79,80,81,82,83,84,86,89,90,91,92,93)

foreach ($r in $rRowsOnly) {
    $wsOikeudet.Range("P$r").Value = "R"
}

# Row 29's P-column value is "R*,W+" -> becomes "R,W"
$wsOikeudet.Range("P29").Value = "R,W"

# Update selection state left behind on the "Roolit" sheet, then return
# focus to "Oikeudet" (which stays the active/visible sheet) with the
# last-edited cell selected.
$wsRoolit.Activate() | Out-Null
$wsRoolit.Range("A6").Select() | Out-Null

$wsOikeudet.Activate() | Out-Null
$wsOikeudet.Range("P93").Select() | Out-Null
